$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("overview_testcases")
$ws.Range("A1").Value = "test"
